$wb = $excel.ActiveWorkbook

# The workbook contains duplicate data tables on the "展览" and "全部类型"
# sheets. This refresh updates the "想去人数" (want-to-go count) column (F)
# for four events on both sheets.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 2311
    $ws.Range("F3").Value = 1763
    $ws.Range("F6").Value = 949
    $ws.Range("F8").Value = 5873
}
